$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "doctor_MA average" column (AF) for rows 4-13 to 1,
# reflecting updated results.
$ws.Range("AF4").Value = 1
$ws.Range("AF5").Value = 1
$ws.Range("AF6").Value = 1
$ws.Range("AF7").Value = 1
$ws.Range("AF8").Value = 1
$ws.Range("AF9").Value = 1
$ws.Range("AF10").Value = 1
$ws.Range("AF11").Value = 1
$ws.Range("AF12").Value = 1
$ws.Range("AF13").Value = 1
